# The graph-description textboxes ("A droite" / "A gauche") on the two
# "Resultats et analyse" slides had their left/right wording swapped
# relative to which plot they actually describe. Fix the wording back
# (commit: "Description des graphes etait inversee").

$p = $ppt.ActivePresentation

# Slide 8 - "F. Resultats et analyse - Avant" - shape "ZoneTexte 8"
$sh8 = $p.Slides.Item(8).Shapes.Item("ZoneTexte 8")
$tr8 = $sh8.TextFrame.TextRange
$tr8.Paragraphs(1, 1).Runs(1, 1).Text = "À gauche, on observe le tracé de la concentration de la solution numérique par rapport à la solution analytique. On voit que ces solutions sont initialement exactes à r = 0.5 mais diverge en s’approchant du centre du cylindre. On peut s’attendre à ce résultat étant donné la condition de Dirichlet qui impose la solution à droite."
$tr8.Paragraphs(3, 1).Runs(1, 1).Text = "À droite, pour ce schéma, on obtient un ordre de convergence / pente d’environ 1 pour les trois erreurs. L’erreur commise est de l’ordre de 10"

# Slide 9 - "F. Resultats et analyse - Centre" - shape "ZoneTexte 8"
$sh9 = $p.Slides.Item(9).Shapes.Item("ZoneTexte 8")
$tr9 = $sh9.TextFrame.TextRange
$tr9.Paragraphs(1, 1).Runs(1, 1).Text = "À gauche, on observe le tracé de la concentration de la solution numérique par rapport à la solution analytique. On voit que ces solutions sont exactement identiques. Cela est conforme aux attentes d’un schéma d’ordre 2, permettant de résoudre exactement une équation différentielle d’ordre 2."
$tr9.Paragraphs(3, 1).Runs(1, 1).Text = "À droite, pour le schéma à discrétisation centrée, on obtient un ordre de convergence / pente non concluant, variant de -6.51 à -3.17. L’erreur observée est de 10"
